$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 39

$ws.Cells.Item($row, 1).Value = 49

# Column B holds a date formatted as plain text ("2026-02-16"); without the
# quote-prefix trick Excel auto-converts it to a date serial number.
$cB = $ws.Cells.Item($row, 2)
$cB.Value = "'2026-02-16"
$cB.Style = "Normal"

$ws.Cells.Item($row, 3).Value = "21:30:24"
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "UP"
$ws.Cells.Item($row, 6).Value = 68771.005

# Columns G and M are empty-string text cells (trade still OPEN, no exit yet).
# A direct Value = "" assignment clears/removes the cell instead of storing an
# empty text value, so use the quote-prefix trick to force a text cell, then
# strip the resulting quote-prefix style.
$cG = $ws.Cells.Item($row, 7)
$cG.Value = "'"
$cG.Style = "Normal"

$ws.Cells.Item($row, 8).Value = "OPEN"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0.75
$ws.Cells.Item($row, 12).Value = "Binance leading with 0.220% move"

$cM = $ws.Cells.Item($row, 13)
$cM.Value = "'"
$cM.Style = "Normal"

$ws.Cells.Item($row, 14).Value = 0
